$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1144
$ws.Range("F5").Value = 59
$ws.Range("F7").Value = 795
$ws.Range("F8").Value = 267
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 429
$ws.Range("F12").Value = 230
$ws.Range("F13").Value = 97
$ws.Range("F14").Value = 947
$ws.Range("F15").Value = 129
$ws.Range("F16").Value = 2058
$ws.Range("F17").Value = 544
$ws.Range("F18").Value = 8953
$ws.Range("F19").Value = 850
$ws.Range("F25").Value = 145

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5612
$ws.Range("F3").Value = 433
$ws.Range("F4").Value = 412

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5612
$ws.Range("F4").Value = 433
$ws.Range("F5").Value = 412
$ws.Range("F7").Value = 1144
$ws.Range("F10").Value = 59
$ws.Range("F12").Value = 795
$ws.Range("F14").Value = 267
$ws.Range("F16").Value = 52
$ws.Range("F17").Value = 429
$ws.Range("F18").Value = 230
$ws.Range("F20").Value = 97
$ws.Range("F22").Value = 947
$ws.Range("F24").Value = 129
$ws.Range("F27").Value = 2058
$ws.Range("F28").Value = 544
$ws.Range("F29").Value = 8953
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 850
